# Apply the "Added 10MHz crystal settings" edit to the MCP2515 Calc workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Calc")

# Crystal Frequency: 16 MHz -> 10 MHz
$ws.Range("B2").Value = 10

# Phase Segment 1 (PS1): 6 -> 2
$ws.Range("G5").Value = 2

# Phase Segment 2 (PS2): 5 -> 3
$ws.Range("G6").Value = 3

# Update the active selection to match the authored state
$ws.Range("G6").Select()

$wb.Save()
